$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44307
$ws.Range("K2").Value = 'Sin especificar'
$ws.Range("M2").Value = 150
$ws.Range("N2").Value = 16000
$ws.Range("O2").Value = 18000
$ws.Range("P2").Value = 17000
$ws.Range("Q2").Value = '$/caja 15 kilos granel'
$ws.Range("R2").Value = "Región de O'Higgins"
$ws.Range("S2").Value = 1133
$ws.Range("T2").Value = 15

# Row 3
$ws.Range("D3").Value = 44280
$ws.Range("K3").Value = 'Sin especificar'
$ws.Range("L3").Value = 'Primera'
$ws.Range("M3").Value = 15
$ws.Range("N3").Value = 360000
$ws.Range("O3").Value = 360000
$ws.Range("P3").Value = 360000
$ws.Range("Q3").Value = '$/bins (450 kilos)'
$ws.Range("S3").Value = 800
$ws.Range("T3").Value = 450

# Row 4
$ws.Range("D4").Value = 44285
$ws.Range("L4").Value = 'Primera'
$ws.Range("M4").Value = 8
$ws.Range("N4").Value = 280000
$ws.Range("O4").Value = 300000
$ws.Range("P4").Value = 290000
$ws.Range("Q4").Value = '$/bins (400 kilos)'
$ws.Range("S4").Value = 725
$ws.Range("T4").Value = 400

# Row 5
$ws.Range("D5").Value = 44266
$ws.Range("L5").Value = 'Segunda'
$ws.Range("M5").Value = 120
$ws.Range("N5").Value = 4800
$ws.Range("O5").Value = 4800
$ws.Range("P5").Value = 4800
$ws.Range("Q5").Value = '$/bandeja 4 kilos'
$ws.Range("R5").Value = 'Provincia del Elquí'
$ws.Range("S5").Value = 1200
$ws.Range("T5").Value = 4

# Row 6
$ws.Range("D6").Value = 44266
$ws.Range("L6").Value = 'Tercera'
$ws.Range("M6").Value = 80
$ws.Range("N6").Value = 4000
$ws.Range("O6").Value = 4000
$ws.Range("P6").Value = 4000
$ws.Range("Q6").Value = '$/bandeja 4 kilos'
$ws.Range("S6").Value = 1000
$ws.Range("T6").Value = 4

# Row 7
$ws.Range("D7").Value = 44312
$ws.Range("K7").Value = 'Wonderfull'
$ws.Range("M7").Value = 24
$ws.Range("N7").Value = 220000
$ws.Range("O7").Value = 240000
$ws.Range("P7").Value = 230000
$ws.Range("Q7").Value = '$/bins (400 kilos)'
$ws.Range("R7").Value = "Región de O'Higgins"
$ws.Range("S7").Value = 575
$ws.Range("T7").Value = 400

# Row 8
$ws.Range("D8").Value = 44312
$ws.Range("M8").Value = 34
$ws.Range("N8").Value = 240000
$ws.Range("O8").Value = 240000
$ws.Range("P8").Value = 240000
$ws.Range("Q8").Value = '$/bins (450 kilos)'
$ws.Range("R8").Value = 'Provincia del Elquí'
$ws.Range("S8").Value = 533
$ws.Range("T8").Value = 450

# Row 9
$ws.Range("D9").Value = 44320
$ws.Range("K9").Value = 'Wonderfull'
$ws.Range("M9").Value = 12
$ws.Range("N9").Value = 250000
$ws.Range("O9").Value = 260000
$ws.Range("P9").Value = 255000
$ws.Range("Q9").Value = '$/bins (400 kilos)'
$ws.Range("R9").Value = 'Provincia de Limarí'
$ws.Range("S9").Value = 638
$ws.Range("T9").Value = 400
